# Auto-generated Excel COM-interop script
# Applies cell value updates to the Adamantoise_Profits price/profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")


# --- ALC ---
$ws1.Range("H101").Value = 747.2222
$ws1.Range("J101").Value = 1124.5
$ws1.Range("L101").Value = 3373.5
$ws1.Range("N101").Value = -6617.5
$ws1.Range("H103").Value = 1031
$ws1.Range("I103").Value = 624.5
$ws1.Range("J103").Value = 1437.5
$ws1.Range("K103").Value = 1873.5
$ws1.Range("L103").Value = 4312.5
$ws1.Range("M103").Value = -1287.5
$ws1.Range("N103").Value = -5484.5
$ws1.Range("H121").Value = 2449.7693
$ws1.Range("J121").Value = 2449.7693
$ws1.Range("L121").Value = 7349.3079
$ws1.Range("N121").Value = -10843.3079
$ws1.Range("H132").Value = 2464.7
$ws1.Range("I132").Value = 2182
$ws1.Range("K132").Value = 6546
$ws1.Range("M132").Value = -4016
$ws1.Range("H135").Value = 1160.16
$ws1.Range("I135").Value = 591.5625
$ws1.Range("J135").Value = 2171
$ws1.Range("K135").Value = 5324.0625
$ws1.Range("L135").Value = 19539
$ws1.Range("M135").Value = -2789.0625
$ws1.Range("N135").Value = -24609
$ws1.Range("H137").Value = 1880996.4
$ws1.Range("I137").Value = 40165.19
$ws1.Range("J137").Value = 6176269
$ws1.Range("K137").Value = 120495.57
$ws1.Range("L137").Value = 18528807
$ws1.Range("M137").Value = -117945.57
$ws1.Range("N137").Value = -18533907
$ws1.Range("H141").Value = 4042.84
$ws1.Range("I141").Value = 4169.8096
$ws1.Range("J141").Value = 3376.25
$ws1.Range("K141").Value = 12509.4288
$ws1.Range("L141").Value = 10128.75
$ws1.Range("M141").Value = -7329.428799999998
$ws1.Range("N141").Value = -20488.75

# --- ARM ---
$ws2.Range("H2").Value = 1774.4783
$ws2.Range("I2").Value = 1107.7693
$ws2.Range("J2").Value = 2641.2
$ws2.Range("K2").Value = 1107.7693
$ws2.Range("L2").Value = 2641.2
$ws2.Range("M2").Value = -994.7692999999999
$ws2.Range("N2").Value = -2867.2
$ws2.Range("H32").Value = 26688364
$ws2.Range("I32").Value = 26915298
$ws2.Range("J32").Value = 23813856
$ws2.Range("K32").Value = 26915298
$ws2.Range("L32").Value = 23813856
$ws2.Range("M32").Value = -26915011
$ws2.Range("N32").Value = -23814430
$ws2.Range("H44").Value = 69984
$ws2.Range("J44").Value = 69984
$ws2.Range("L44").Value = 69984
$ws2.Range("N44").Value = -70960
$ws2.Range("H74").Value = 1567.1333
$ws2.Range("I74").Value = 1448.9412
$ws2.Range("K74").Value = 1448.9412
$ws2.Range("M74").Value = -574.9412
$ws2.Range("H77").Value = 1567.1333
$ws2.Range("I77").Value = 1448.9412
$ws2.Range("K77").Value = 7244.706
$ws2.Range("M77").Value = -2876.706
$ws2.Range("H116").Value = 1774.4783
$ws2.Range("I116").Value = 1107.7693
$ws2.Range("J116").Value = 2641.2
$ws2.Range("K116").Value = 1107.7693
$ws2.Range("L116").Value = 2641.2
$ws2.Range("M116").Value = 1186.2307
$ws2.Range("N116").Value = -7229.2

# --- BSM ---
$ws3.Range("H3").Value = 1774.4783
$ws3.Range("I3").Value = 1107.7693
$ws3.Range("J3").Value = 2641.2
$ws3.Range("K3").Value = 1107.7693
$ws3.Range("L3").Value = 2641.2
$ws3.Range("M3").Value = -993.7692999999999
$ws3.Range("N3").Value = -2869.2
$ws3.Range("H11").Value = 859
$ws3.Range("J11").Value = 800
$ws3.Range("L11").Value = 800
$ws3.Range("N11").Value = -1080
$ws3.Range("H20").Value = 67100.5
$ws3.Range("I20").Value = 103761
$ws3.Range("K20").Value = 103761
$ws3.Range("M20").Value = -103514
$ws3.Range("H22").Value = 249
$ws3.Range("I22").Value = 249
$ws3.Range("K22").Value = 249
$ws3.Range("M22").Value = -76
$ws3.Range("H94").Value = 598.625
$ws3.Range("I94").Value = 554.0909
$ws3.Range("K94").Value = 554.0909
$ws3.Range("M94").Value = -103.0909
$ws3.Range("H99").Value = 2175.7646
$ws3.Range("I99").Value = 1540.6666
$ws3.Range("K99").Value = 1540.6666
$ws3.Range("M99").Value = -42.66660000000002
$ws3.Range("H107").Value = 1586.875
$ws3.Range("I107").Value = 1242.1428
$ws3.Range("K107").Value = 1242.1428
$ws3.Range("M107").Value = 677.8571999999999
$ws3.Range("H134").Value = 4765003.5
$ws3.Range("I134").Value = 5558462
$ws3.Range("J134").Value = 4252
$ws3.Range("K134").Value = 16675386
$ws3.Range("L134").Value = 12756
$ws3.Range("M134").Value = -16672851
$ws3.Range("N134").Value = -17826

# --- CRP ---
$ws4.Range("H107").Value = 1059.4
$ws4.Range("I107").Value = 474.33334
$ws4.Range("K107").Value = 474.33334
$ws4.Range("M107").Value = 1445.66666
$ws4.Range("H134").Value = 2004.2858
$ws4.Range("I134").Value = 1520.125
$ws4.Range("K134").Value = 4560.375
$ws4.Range("M134").Value = -2025.375

# --- CUL ---
$ws5.Range("H5").Value = 1456.6
$ws5.Range("I5").Value = 648.5
$ws5.Range("K5").Value = 1945.5
$ws5.Range("M5").Value = -1833.5
$ws5.Range("H14").Value = 10366.182
$ws5.Range("I14").Value = 10366.182
$ws5.Range("K14").Value = 31098.546
$ws5.Range("M14").Value = -30925.546
$ws5.Range("H64").Value = 0
$ws5.Range("I64").Value = 0
$ws5.Range("J64").Value = 0
$ws5.Range("K64").Value = 0
$ws5.Range("L64").Value = 0
$ws5.Range("M64").ClearContents()
$ws5.Range("N64").ClearContents()
$ws5.Range("H67").Value = 0
$ws5.Range("I67").Value = 0
$ws5.Range("J67").Value = 0
$ws5.Range("K67").Value = 0
$ws5.Range("L67").Value = 0
$ws5.Range("M67").ClearContents()
$ws5.Range("N67").ClearContents()
$ws5.Range("H68").Value = 856.4286
$ws5.Range("J68").Value = 865.5
$ws5.Range("L68").Value = 2596.5
$ws5.Range("N68").Value = -4218.5
$ws5.Range("H71").Value = 856.4286
$ws5.Range("J71").Value = 865.5
$ws5.Range("L71").Value = 7789.5
$ws5.Range("N71").Value = -15901.5
$ws5.Range("H80").Value = 2999
$ws5.Range("J80").Value = 2999
$ws5.Range("L80").Value = 8997
$ws5.Range("N80").Value = -10869
$ws5.Range("H83").Value = 2999
$ws5.Range("J83").Value = 2999
$ws5.Range("L83").Value = 26991
$ws5.Range("N83").Value = -36351
$ws5.Range("H97").Value = 266.33334
$ws5.Range("I97").Value = 600
$ws5.Range("J97").Value = 99.5
$ws5.Range("K97").Value = 1800
$ws5.Range("L97").Value = 298.5
$ws5.Range("M97").Value = -1304
$ws5.Range("N97").Value = -1290.5
$ws5.Range("H107").Value = 868.06665
$ws5.Range("I107").Value = 772.6667
$ws5.Range("K107").Value = 2318.0001
$ws5.Range("M107").Value = -398.0001000000002
$ws5.Range("H119").Value = 800
$ws5.Range("I119").Value = 800
$ws5.Range("K119").Value = 2400
$ws5.Range("M119").Value = 2438
$ws5.Range("H135").Value = 1456.6
$ws5.Range("I135").Value = 648.5
$ws5.Range("K135").Value = 5836.5
$ws5.Range("M135").Value = -3301.5

# --- GSM ---
$ws6.Range("H41").Value = 1370.5
$ws6.Range("I41").Value = 1370.5
$ws6.Range("K41").Value = 1370.5
$ws6.Range("M41").Value = -1015.5
$ws6.Range("H62").Value = 34000
$ws6.Range("I62").Value = 34000
$ws6.Range("K62").Value = 34000
$ws6.Range("M62").Value = -33314
$ws6.Range("H65").Value = 34000
$ws6.Range("I65").Value = 34000
$ws6.Range("K65").Value = 102000
$ws6.Range("M65").Value = -98568
$ws6.Range("H70").Value = 17872.137
$ws6.Range("I70").Value = 102646.664
$ws6.Range("J70").Value = 4486.684
$ws6.Range("K70").Value = 102646.664
$ws6.Range("L70").Value = 4486.684
$ws6.Range("M70").Value = -102376.664
$ws6.Range("N70").Value = -5026.684
$ws6.Range("H73").Value = 17872.137
$ws6.Range("I73").Value = 102646.664
$ws6.Range("J73").Value = 4486.684
$ws6.Range("K73").Value = 102646.664
$ws6.Range("L73").Value = 4486.684
$ws6.Range("M73").Value = -101710.664
$ws6.Range("N73").Value = -6358.684
$ws6.Range("H97").Value = 519.3929000000001
$ws6.Range("I97").Value = 401.28
$ws6.Range("K97").Value = 401.28
$ws6.Range("M97").Value = 94.72000000000003

# --- LTW ---
$ws7.Range("H132").Value = 716733.6
$ws7.Range("I132").Value = 1252298.8
$ws7.Range("K132").Value = 3756896.4
$ws7.Range("M132").Value = -3754366.4

# --- WVR ---
$ws8.Range("H70").Value = 37463
$ws8.Range("J70").Value = 37463
$ws8.Range("L70").Value = 37463
$ws8.Range("N70").Value = -38093
$ws8.Range("H73").Value = 37463
$ws8.Range("J73").Value = 37463
$ws8.Range("L73").Value = 37463
$ws8.Range("N73").Value = -39647
$ws8.Range("H132").Value = 58742.61
$ws8.Range("I132").Value = 69510.8
$ws8.Range("J132").Value = 4901.6665
$ws8.Range("K132").Value = 208532.4
$ws8.Range("L132").Value = 14704.9995
$ws8.Range("M132").Value = -206002.4
$ws8.Range("N132").Value = -19764.9995
$ws8.Range("H136").Value = 1182.8182
$ws8.Range("I136").Value = 1182.8182
$ws8.Range("K136").Value = 3548.4546
$ws8.Range("M136").Value = -998.4546
